# Weekly update: add a new "Poroto granado" price observation as a new
# row right before the current row 37, pushing all later rows down by one
# (matches the commit "Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 37; everything from old row 37
# downward shifts to row+1 (old 37 -> 38, ..., old 57 -> 58).
$ws.Rows.Item(37).Insert()

# Populate the new row 37 with the new weekly data point.
$ws.Cells.Item(37, 1).Value  = 2
$ws.Cells.Item(37, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(37, 3).Value  = "Coquimbo"
$ws.Cells.Item(37, 4).Value  = 44553
$ws.Cells.Item(37, 5).Value  = 4
$ws.Cells.Item(37, 6).Value  = 100112030
$ws.Cells.Item(37, 7).Value  = "Poroto granado"
$ws.Cells.Item(37, 8).Value  = "Sin especificar"
$ws.Cells.Item(37, 9).Value  = "Primera"
$ws.Cells.Item(37, 10).Value = 600
$ws.Cells.Item(37, 11).Value = 25000
$ws.Cells.Item(37, 12).Value = 27000
$ws.Cells.Item(37, 13).Value = 26000
$ws.Cells.Item(37, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(37, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(37, 16).Value = 1040
$ws.Cells.Item(37, 17).Value = 25
$ws.Cells.Item(37, 18).Value = "Hortaliza"
